$d = $word.ActiveDocument

# Update the date/day heading at the top of the worksheet.
$d.Content.Find.Execute("2023-11-03 Friday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2023-11-04 Saturday", 2)

# Update every division problem in the practice table. We address cells
# directly by (row, column) instead of doing a blind global find/replace
# because one of the new values ("81÷9=") collides with an old value
# used elsewhere in the table, which a naive text-replace pass could
# clobber depending on execution order.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "71÷4="
$t.Cell(1, 2).Range.Text = "70÷8="
$t.Cell(1, 3).Range.Text = "72÷4="
$t.Cell(1, 4).Range.Text = "36÷9="
$t.Cell(1, 5).Range.Text = "52÷7="

$t.Cell(5, 1).Range.Text = "38÷9="
$t.Cell(5, 2).Range.Text = "47÷8="
$t.Cell(5, 3).Range.Text = "81÷9="
$t.Cell(5, 4).Range.Text = "29÷9="
$t.Cell(5, 5).Range.Text = "86÷6="

$t.Cell(9, 1).Range.Text = "60÷9="
$t.Cell(9, 2).Range.Text = "99÷3="
$t.Cell(9, 3).Range.Text = "55÷7="
$t.Cell(9, 4).Range.Text = "54÷4="
$t.Cell(9, 5).Range.Text = "30÷2="

$t.Cell(13, 1).Range.Text = "48÷7="
$t.Cell(13, 2).Range.Text = "80÷2="
$t.Cell(13, 3).Range.Text = "38÷7="
$t.Cell(13, 4).Range.Text = "12÷3="
$t.Cell(13, 5).Range.Text = "15÷8="

$t.Cell(17, 1).Range.Text = "30÷7="
$t.Cell(17, 2).Range.Text = "67÷7="
$t.Cell(17, 3).Range.Text = "17÷8="
$t.Cell(17, 4).Range.Text = "32÷2="
$t.Cell(17, 5).Range.Text = "39÷6="
